# "this verison is working"
#
# Changes applied:
#  1. Reorder worksheets so "swathi" comes before "raj" (swathi becomes the
#     first tab, raj the second - raj stays the active tab).
#  2. Remove the leading "20105/" zip code from the "zipcodes" value (cell
#     D2) on both the "raj" and "swathi" sheets.
#  3. Add a new "Backup" worksheet at the end that keeps a copy of the
#     original (pre-edit), full zip code list in cell A1.
#  4. Leave behind the selection state that was captured when the file was
#     last saved (swathi!D2 selected, raj!E5 selected, raj active).

$wb = $excel.ActiveWorkbook

$raj = $wb.Worksheets.Item("raj")
$swathi = $wb.Worksheets.Item("swathi")

# --- 1. Reorder sheets: swathi, then raj -------------------------------
$swathi.Move($raj)

# NOTE: worksheet object references returned by the COM model are bound to
# sheet *position*, not identity, so after a Move() the old variables now
# point at whatever sheet occupies that slot.  Re-fetch by name.
$raj = $wb.Worksheets.Item("raj")
$swathi = $wb.Worksheets.Item("swathi")

# --- 2. Update the zip code list on both sheets -------------------------
$oldZips = "20105/27560/75024/36101/99801/85001/72201/94203/80201/06101/19901/20001/32301/30301/96801/83701/62701/46201/50301/66601/40601/70801/04330/21401/02108/48901/55101/39201/65101/59601/68501/89701/03301/08601/87501/12201/27601/58501/43201/73101/93701/17101/02901/29201/37201/73301/84101/05601/23218/98501/25301"
$newZips = "27560/75024/36101/99801/85001/72201/94203/80201/06101/19901/20001/32301/30301/96801/83701/62701/46201/50301/66601/40601/70801/04330/21401/02108/48901/55101/39201/65101/59601/68501/89701/03301/08601/87501/12201/27601/58501/43201/73101/93701/17101/02901/29201/37201/73301/84101/05601/23218/98501/25301"

$raj.Range("D2").Value = $newZips
$swathi.Range("D2").Value = $newZips

# --- 3. Add the Backup sheet with the original zip code list ------------
$backup = $wb.Worksheets.Add()
$backup.Name = "Backup"
$backup.Range("A1").Value = $oldZips

# Move Backup to the end of the sheet tabs.
$backup.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Re-fetch again since positions shifted after the Backup move.
$raj = $wb.Worksheets.Item("raj")
$swathi = $wb.Worksheets.Item("swathi")

# --- 4. Restore selections / active sheet --------------------------------
$swathi.Range("D2").Select()
$raj.Range("E5").Select()
$raj.Activate()
